$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Range("D9").Hyperlinks.Add($ws.Range("D9"), "mailto:edu_21_91@hotmail.com")
$ws.Range("D9").Style = "Normal"
